$wb = $excel.ActiveWorkbook

$oldGuid = "0f44ba9f-98f1-4225-9f7e-b97f711a48a0"
$newGuid = "9580bc68-c15d-4358-a815-881302bafc9b"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}
$wsOverview.Range("G2").Value = "2016-08-26 02:57:58"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
$wsZhCn.Range("G2").Value = "$newGuid.6c79225986249f9c509165f64bd220448001178e.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 02:57:52"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
$wsDeDe.Range("G2").Value = "$newGuid.6c79225986249f9c509165f64bd220448001178e.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 02:57:58"
